# feat: generate exam without images
#
# Remove the rows that represent the "with images" variant of each
# exercise (these rows have no value in the "title" (D) column and no
# value in the "image bottom description" (G) column) so that only the
# text-only ("without images") rows remain.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

# Collect the row numbers (1-based, including header row 1) whose D and G
# cells are both empty - these are the rows to delete.
$rowsToDelete = @()
for ($r = 2; $r -le $lastRow; $r++) {
    $dVal = $ws.Cells.Item($r, 4).Value2
    $gVal = $ws.Cells.Item($r, 7).Value2
    $dEmpty = ($dVal -eq $null) -or ($dVal -eq "")
    $gEmpty = ($gVal -eq $null) -or ($gVal -eq "")
    if ($dEmpty -and $gEmpty) {
        $rowsToDelete += $r
    }
}

# Delete rows from bottom to top so row numbers of not-yet-processed rows
# don't shift while we iterate.
for ($i = $rowsToDelete.Count - 1; $i -ge 0; $i--) {
    $ws.Rows.Item($rowsToDelete[$i]).Delete()
}
